$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill B3:C22 with explicit 0 values and clear their formatting (remove fill/alignment)
$rng = $ws.Range("B3:C22")
$rng.Value = 0
$rng.ClearFormats()

# 2. Update summary formulas in row 23 and 24 to cover the full B3:B22 / C3:C22 ranges
$ws.Range("E23").Formula = "=SUM(E3:E22)/SUM(B3:B22)"
$ws.Range("E24").Formula = "=SUM(F3:F22)/SUM(C3:C22)"

# 3. Apply center alignment to F23 (new formatted, empty cell)
$ws.Range("F23").HorizontalAlignment = -4108

# 4. Move the active selection to F27
$ws.Range("F27").Select()
